$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 6-9 (only 5 rows of data remain: header + 4 data rows)
$ws.Range("A6:T9").EntireRow.Delete()

# Row 2: FAPs, Lgi3, Adam23, ECs
$ws.Range("B2").Value = "Lgi3"
$ws.Range("C2").Value = "Adam23"
$ws.Range("D2").Value = "ECs"
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("M2").Value = 0.165747
$ws.Range("N2").Value = 0.497241
$ws.Range("O2").Value = 0.008095785894995438
$ws.Range("P2").Value = 0.00809578589499544
$ws.Range("Q2").Value = 0.254549435937
$ws.Range("R2").Value = 2.290944923433
$ws.Range("S2").Value = 0.008095785894995438
$ws.Range("T2").Value = 0.00809578589499544

# Row 3: FAPs, Lgi3, Adam23, FAPs
$ws.Range("B3").Value = "Lgi3"
$ws.Range("C3").Value = "Adam23"
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("O3").Value = 0.7079722685862583
$ws.Range("P3").Value = 0.7079722685862583
$ws.Range("S3").Value = 0.7079722685862583
$ws.Range("T3").Value = 0.7079722685862583

# Row 4: FAPs, Lgi3, Adam23, MuSCs
$ws.Range("B4").Value = "Lgi3"
$ws.Range("C4").Value = "Adam23"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("M4").Value = 5.642879333333333
$ws.Range("N4").Value = 16.928638
$ws.Range("O4").Value = 0.2756221404547972
$ws.Range("P4").Value = 0.2756221404547972
$ws.Range("Q4").Value = 8.666170436632667
$ws.Range("R4").Value = 77.99553392969401
$ws.Range("S4").Value = 0.2756221404547972
$ws.Range("T4").Value = 0.2756221404547972

# Row 5: FAPs, Lgi3, Adam23, Resolving-Mac
$ws.Range("B5").Value = "Lgi3"
$ws.Range("C5").Value = "Adam23"
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.1701286666666667
$ws.Range("N5").Value = 0.510386
$ws.Range("O5").Value = 0.008309805063949155
$ws.Range("P5").Value = 0.008309805063949155
$ws.Range("Q5").Value = 0.2612786725353334
$ws.Range("R5").Value = 2.351508052818
$ws.Range("S5").Value = 0.008309805063949155
$ws.Range("T5").Value = 0.008309805063949155
